$d = $word.ActiveDocument

# --- Step 1 -----------------------------------------------------------
# Paragraph: "The battery installed in the lamp is a single cell 400mAh
# li-ion battery ... external short damages. The user can put the lamp
# into two modes when charging the battery:"
# followed by a bulleted paragraph beginning "Display charging mode: the
# user needs to plug ..."
#
# Collapse the two modes sentence and the bold "Display charging mode:"
# label away, merging the following list paragraph back into the body
# paragraph so the text reads straight on to "The user needs to plug the
# USB end ...". (^p matches the paragraph mark that separates the two
# paragraphs.)
$r1 = $d.Content
$found1 = $r1.Find.Execute( `
    "he user can put the lamp into two modes when charging the battery:^pDisplay charging mode: t", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
if (-not $found1) {
    throw "Could not find the 'two modes / Display charging mode' text to replace"
}

# --- Step 2 -----------------------------------------------------------
# Remove the whole "Dark charging mode: ..." bulleted paragraph, merging
# its neighbours (the paragraph mark before it and the one after it) so
# the narrative continues directly into "A typical charge cycle takes
# two hours ...".
$r2 = $d.Content
$found2 = $r2.Find.Execute( `
    "^pDark charging mode: the user needs to plug the round end of the charging cable into the bottom of the lamp first, then plug the USB end to the charger, no light should come up during a charging cycle. ^p", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 2)
if (-not $found2) {
    throw "Could not find the 'Dark charging mode' paragraph to remove"
}

# --- Step 3 -----------------------------------------------------------
# Word's "_GoBack" bookmark marks the site of the author's last edit.
# Relocate it from right after "installed" to right before "A typical
# charge cycle ...", matching where the editing session left off.
$anchor = $d.Content
$foundAnchor = $anchor.Find.Execute("A typical charge cycle", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $foundAnchor) {
    throw "Could not find the 'A typical charge cycle' anchor text"
}

$goBackRange = $d.Range($anchor.Start, $anchor.Start)
try {
    $existing = $d.Bookmarks("_GoBack")
    $existing.Delete()
} catch {
    # no pre-existing _GoBack bookmark - nothing to remove
}
$d.Bookmarks.Add("_GoBack", $goBackRange)
